$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "2025-09-25"
$ws.Range("A41").Style = "Normal"
$ws.Range("B41").Value = 56.63999938964844
$ws.Range("C41").Value = 664.2999877929688
$ws.Range("D41").Value = 332.25
